# Auto-generated edit script: updates numeric cells (columns H-N)
# on each leve-profit sheet to match the scheduled-runner recompute.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 150.93103
$ws.Range("I33").Value = 152.28572
$ws.Range("K33").Value = 152.28572
$ws.Range("M33").Value = 76.71428
$ws.Range("H43").Value = 3625
$ws.Range("J43").Value = 4166.6665
$ws.Range("L43").Value = 4166.6665
$ws.Range("N43").Value = -4304.6665
$ws.Range("H96").Value = 2167.3333
$ws.Range("I96").Value = 706.25
$ws.Range("J96").Value = 5089.5
$ws.Range("K96").Value = 2118.75
$ws.Range("L96").Value = 15268.5
$ws.Range("M96").Value = -745.75
$ws.Range("N96").Value = -18014.5
$ws.Range("H101").Value = 1436.3889
$ws.Range("I101").Value = 1171.7
$ws.Range("J101").Value = 1767.25
$ws.Range("K101").Value = 3515.1
$ws.Range("L101").Value = 5301.75
$ws.Range("M101").Value = -1893.1
$ws.Range("N101").Value = -8545.75
$ws.Range("H137").Value = 14935143
$ws.Range("I137").Value = 720782.5600000001
$ws.Range("J137").Value = 37046372
$ws.Range("K137").Value = 2162347.68
$ws.Range("L137").Value = 111139116
$ws.Range("M137").Value = -2159797.68
$ws.Range("N137").Value = -111144216
$ws.Range("H138").Value = 4175.8643
$ws.Range("I138").Value = 1268.7241
$ws.Range("K138").Value = 3806.1723
$ws.Range("M138").Value = 1333.8277

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1784.9
$ws.Range("I32").Value = 1779.707
$ws.Range("J32").Value = 2299
$ws.Range("K32").Value = 1779.707
$ws.Range("L32").Value = 2299
$ws.Range("M32").Value = -1492.707
$ws.Range("N32").Value = -2873
$ws.Range("H45").Value = 2879.9
$ws.Range("I45").Value = 2575.923
$ws.Range("K45").Value = 2575.923
$ws.Range("M45").Value = -2198.923
$ws.Range("H61").Value = 4292.9287
$ws.Range("I61").Value = 4096.7437
$ws.Range("J61").Value = 4743
$ws.Range("K61").Value = 4096.7437
$ws.Range("L61").Value = 4743
$ws.Range("M61").Value = -3884.7437
$ws.Range("N61").Value = -5167
$ws.Range("H74").Value = 1341
$ws.Range("I74").Value = 1295.6
$ws.Range("J74").Value = 1386.4
$ws.Range("K74").Value = 1295.6
$ws.Range("L74").Value = 1386.4
$ws.Range("M74").Value = -421.5999999999999
$ws.Range("N74").Value = -3134.4
$ws.Range("H77").Value = 1341
$ws.Range("I77").Value = 1295.6
$ws.Range("J77").Value = 1386.4
$ws.Range("K77").Value = 6478
$ws.Range("L77").Value = 6932
$ws.Range("M77").Value = -2110
$ws.Range("N77").Value = -15668
$ws.Range("H102").Value = 1784.0625
$ws.Range("I102").Value = 1830.3846
$ws.Range("K102").Value = 1830.3846
$ws.Range("M102").Value = -208.3846000000001
$ws.Range("H122").Value = 5425.2
$ws.Range("I122").Value = 3087.8
$ws.Range("K122").Value = 9263.400000000001
$ws.Range("M122").Value = -6813.400000000001
$ws.Range("H136").Value = 4292.9287
$ws.Range("I136").Value = 4096.7437
$ws.Range("J136").Value = 4743
$ws.Range("K136").Value = 12290.2311
$ws.Range("L136").Value = 14229
$ws.Range("M136").Value = -9740.231100000001
$ws.Range("N136").Value = -19329

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("N9").Value = $null
$ws.Range("H44").Value = 27886.428
$ws.Range("J44").Value = 27619.4
$ws.Range("L44").Value = 27619.4
$ws.Range("N44").Value = -28613.4
$ws.Range("H105").Value = 2712.6667
$ws.Range("I105").Value = 2709.9048
$ws.Range("J105").Value = 2732
$ws.Range("K105").Value = 2709.9048
$ws.Range("L105").Value = 2732
$ws.Range("M105").Value = -962.9047999999998
$ws.Range("N105").Value = -6226
$ws.Range("H116").Value = 76451.73
$ws.Range("J116").Value = 76451.73
$ws.Range("L116").Value = 76451.73
$ws.Range("N116").Value = -85629.73
$ws.Range("H134").Value = 5094.0166
$ws.Range("I134").Value = 2479.5134
$ws.Range("J134").Value = 9299.956
$ws.Range("K134").Value = 7438.540199999999
$ws.Range("L134").Value = 27899.868
$ws.Range("M134").Value = -4903.540199999999
$ws.Range("N134").Value = -32969.868

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2227
$ws.Range("I58").Value = 2227.4
$ws.Range("K58").Value = 2227.4
$ws.Range("M58").Value = -2024.4
$ws.Range("H107").Value = 725.8461
$ws.Range("I107").Value = 678
$ws.Range("K107").Value = 678
$ws.Range("M107").Value = 1242
$ws.Range("H136").Value = 2227
$ws.Range("I136").Value = 2227.4
$ws.Range("K136").Value = 6682.200000000001
$ws.Range("M136").Value = -4132.200000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1062.1875
$ws.Range("I5").Value = 258.2857
$ws.Range("J5").Value = 1687.4445
$ws.Range("K5").Value = 774.8571000000001
$ws.Range("L5").Value = 5062.333500000001
$ws.Range("M5").Value = -662.8571000000001
$ws.Range("N5").Value = -5286.333500000001
$ws.Range("H135").Value = 1062.1875
$ws.Range("I135").Value = 258.2857
$ws.Range("J135").Value = 1687.4445
$ws.Range("K135").Value = 2324.5713
$ws.Range("L135").Value = 15187.0005
$ws.Range("M135").Value = 210.4286999999999
$ws.Range("N135").Value = -20257.0005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").Value = $null
$ws.Range("H93").Value = 44900.668
$ws.Range("J93").Value = 44900.668
$ws.Range("L93").Value = 44900.668
$ws.Range("N93").Value = -48644.668
$ws.Range("H122").Value = 349639.72
$ws.Range("J122").Value = 7499.143
$ws.Range("L122").Value = 22497.429
$ws.Range("N122").Value = -27397.429

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 934.85
$ws.Range("I22").Value = 834.2308
$ws.Range("J22").Value = 1121.7142
$ws.Range("K22").Value = 834.2308
$ws.Range("L22").Value = 1121.7142
$ws.Range("M22").Value = -539.2308
$ws.Range("N22").Value = -1711.7142
$ws.Range("H27").Value = 934.85
$ws.Range("I27").Value = 834.2308
$ws.Range("J27").Value = 1121.7142
$ws.Range("K27").Value = 834.2308
$ws.Range("L27").Value = 1121.7142
$ws.Range("M27").Value = -727.2308
$ws.Range("N27").Value = -1335.7142
$ws.Range("H48").Value = 47946
$ws.Range("J48").Value = 47946
$ws.Range("L48").Value = 47946
$ws.Range("N48").Value = -49268
$ws.Range("H61").Value = 2196.6
$ws.Range("I61").Value = 1981.1428
$ws.Range("K61").Value = 1981.1428
$ws.Range("M61").Value = -1779.1428
$ws.Range("H93").Value = 3535.6667
$ws.Range("I93").Value = 2399
$ws.Range("K93").Value = 2399
$ws.Range("M93").Value = -1151
$ws.Range("H97").Value = 59169.5
$ws.Range("J97").Value = 59169.5
$ws.Range("L97").Value = 59169.5
$ws.Range("N97").Value = -61151.5
$ws.Range("H100").Value = 2295
$ws.Range("I100").Value = 2354.8
$ws.Range("K100").Value = 2354.8
$ws.Range("M100").Value = -1813.8
$ws.Range("H113").Value = 2196.6
$ws.Range("I113").Value = 1981.1428
$ws.Range("K113").Value = 1981.1428
$ws.Range("M113").Value = 188.8571999999999
$ws.Range("H116").Value = 120000
$ws.Range("J116").Value = 120000
$ws.Range("L116").Value = 120000
$ws.Range("N116").Value = -129178
$ws.Range("H122").Value = 7695.4136
$ws.Range("I122").Value = 3848.7334
$ws.Range("J122").Value = 11816.857
$ws.Range("K122").Value = 11546.2002
$ws.Range("L122").Value = 35450.571
$ws.Range("M122").Value = -9096.200199999999
$ws.Range("N122").Value = -40350.571
$ws.Range("H136").Value = 4992.4585
$ws.Range("I136").Value = 3059.1333
$ws.Range("K136").Value = 9177.3999
$ws.Range("M136").Value = -6627.3999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3508.606
$ws.Range("I122").Value = 2754.158
$ws.Range("J122").Value = 4532.5
$ws.Range("K122").Value = 8262.474
$ws.Range("L122").Value = 13597.5
$ws.Range("M122").Value = -5812.474
$ws.Range("N122").Value = -18497.5

Write-Host "Applied 204 cell updates across 8 sheets"
